# "blog dan comment belom!!"
#
# Re-label the "kategori" column (N) placeholder values from the old
# "Nyoba aja / Nyoba doang / Nyoba atuh / Nyoba yaaa" set down to a
# smaller "nyoba_aja / nyoba_doang / nyoba_1" set, and move the scroll
# position / active selection down to the bottom of the table
# (topLeftCell A48, active cell M54) like the author left it before
# committing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- kategori (column N) relabeling ------------------------------------
# rows 2-19  -> "nyoba_aja"
# rows 20-48 -> "nyoba_doang"
# row 49     -> "nyoba_1"
# (the old strings "Nyoba aja"/"Nyoba doang"/"Nyoba atuh"/"Nyoba yaaa" are
# thereby fully dereferenced and dropped from the shared-string table; the
# "foto" placeholder-image URL used by columns O:S is untouched and simply
# shifts down in the table once the unused entries ahead of it are gone)
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 14).Value = "nyoba_aja"
}
for ($r = 20; $r -le 48; $r++) {
    $ws.Cells.Item($r, 14).Value = "nyoba_doang"
}
$ws.Cells.Item(49, 14).Value = "nyoba_1"

# --- scroll position / selection ---------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 48
$win.ScrollColumn = 1
$ws.Range("M54").Select()
